$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Rename the metadata4Ing headers to metadata4ing (lowercase "ing")
$ws.Range("D1").Value = "metadata4ing_IRI"
$ws.Range("E1").Value = "metadata4ing_DESC"

# Add the new F column header, matching the header style used by B1:E1
$ws.Range("F1").Value = "metadata4ing_DEF"
$ws.Range("F1").Font.Bold = $true
$ws.Range("F1").HorizontalAlignment = -4108
$ws.Range("F1").VerticalAlignment = -4160
$ws.Range("F1").Borders.LineStyle = 1

# Add the new F column values (definitions) for each data row
$ws.Range("F2").Value = "['p is a process if p is an occurrent that has temporal proper parts and for some time t, p specifically depends on some material entity at t. [BFO]', locstr(" + [char]34 + "Process, i.e., a physical entity with a temporal evolution that 'has a meaning for the ontologist'" + [char]34 + ", 'en')]"
$ws.Range("F3").Value = "['To say that b is a realizable entity is to say that b is a specifically dependent continuant that inheres in some independent continuant which is not a spatial region and is of a type instances of which are realized in processes of a correlated type.´[BFO]', 'To say that b is a realizable entity is to say that b is a specifically dependent continuant that inheres in some independent continuant which is not a spatial region and is of a type instances of which are realized in processes of a correlated type. (axiom label in BFO2 Reference: [058-002])']"
$ws.Range("F4").Value = "[locstr('A role is the function of an entity or agent with respect to an activity, in the context of a usage, generation, invalidation, association, start, and end.', 'en')]"
